# Insert a new "is_targeted list" sheet (with TRUE/FALSE values) right after
# "analyte_class list", matching the new workbook.xml sheet ordering produced by
# the commit, and point the N column's data validation at it instead of the
# old inline "TRUE,FALSE" literal-list formula.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("analyte_class list")
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $afterSheet)
$newSheet.Name = "is_targeted list"

# Write literal text "TRUE" / "FALSE" (not boolean cells) into A1:A2, matching
# the other "* list" helper sheets which store their options as plain text.
$newSheet.Cells.Item(1, 1).Value = "'TRUE"
$newSheet.Cells.Item(2, 1).Value = "'FALSE"
$newSheet.Range("A1:A2").Style = "Normal"

# Repoint the is_targeted (column N) data validation at the new list sheet,
# in place, so it keeps its original position among the sheet's validations.
$mainSheet = $wb.Worksheets.Item("Export as TSV")
$targetRange = $mainSheet.Range("N2:N1048576")
$targetRange.Validation.Modify(3, 1, 1, "='is_targeted list'!`$A`$1:`$A`$2")
$targetRange.Validation.ErrorTitle = "Value must come from list"
$targetRange.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."

# Restore the original active sheet/selection (inserting a sheet makes it
# active by default) so we don't introduce an unrelated view-state change.
$mainSheet.Activate()

Write-Output "Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
